# Apply crypto price/volume/name updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.762.64'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = '2.479.57'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.79%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.554'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.84%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.519'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.27'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E11').Value = '  +9.41%  '
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').Value = '2.860.96'
$ws.Range('E13').Value = '  +0.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.91'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.80'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range('D16').Value = '2.469.31'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.793'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.99%  '
$ws.Range('D18').Value = '41.714.08'
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').Value = '0.0₃0952'
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.29'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.34'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.28'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.78%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.76'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.17%  '
$ws.Range('E25').Value = '  +2.49%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.81'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('E28').Value = '  +1.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.84'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.21'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.77'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.01%  '
$ws.Range('E32').Value = '  +1.67%  '
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0771'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.44%  '
$ws.Range('E35').Value = '  +0.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.38'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.71%  '
$ws.Range('E37').Value = '  +5.05%  '
$ws.Range('E38').Value = '  +1.81%  '
$ws.Range('E39').Value = '  +1.76%  '
$ws.Range('E40').Value = '  +0.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.04'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('E42').Value = '  +9.55%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.001.58'
$ws.Range('E43').Value = '  +2.97%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.50'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.53%  '
$ws.Range('E45').Value = '  +1.08%  '
$ws.Range('E46').Value = '  +1.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.35'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.17%  '
$ws.Range('D48').Value = '2.716.11'
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '97.67'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.46'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '67.33'
$ws.Range('D51').Style = 'Normal'
